$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected reference values (column A = fmod, column B = R) for rows 3-19.
# The previous last row (20) is removed entirely as part of the correction.
$A = @(
    "11.439509",
    "27.74356",
    "30.087986",
    "32.038048",
    "35.0145",
    "37.96706",
    "41.273243",
    "43.594078",
    "45.3298",
    "48.881107",
    "53.121265",
    "54.943813",
    "56.827637",
    "59.387787",
    "64.01713",
    "73.04394",
    "157.47661"
)

$B = @(
    "0.102054834",
    "0.3669846",
    "0.39998206",
    "0.426939",
    "0.45447007",
    "0.47624823",
    "0.49128193",
    "0.49116743",
    "0.48724157",
    "0.47324988",
    "0.45366368",
    "0.43956152",
    "0.42367226",
    "0.40197152",
    "0.35809588",
    "0.28570408",
    "0.077257395"
)

# Write column A first, then column B, so the values are stored as text
# (matching the original file's string-typed cells) while keeping the
# default "Normal" cell style (no explicit number format override).
for ($i = 0; $i -lt $A.Length; $i++) {
    $row = 3 + $i
    $cA = $ws.Cells.Item($row, 1)
    $cA.NumberFormat = "@"
    $cA.Value = $A[$i]
    $cA.Style = "Normal"
}

for ($i = 0; $i -lt $B.Length; $i++) {
    $row = 3 + $i
    $cB = $ws.Cells.Item($row, 2)
    $cB.NumberFormat = "@"
    $cB.Value = $B[$i]
    $cB.Style = "Normal"
}

# Remove the now-obsolete last row (previously row 20) entirely.
$ws.Rows.Item(20).Delete()

$ws.Range("B3:B19").Select()
